$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:F (error metrics) for rows 2-11 shift down by one row: row N's old
# values move to row N+1, the old row 11 values are dropped, and a new set of
# (near-zero) values is inserted at row 2. Column G is simply incremented by 1
# for every row (2-11).

$newRow2 = @(0.0000002026185291058083, 0.0000006752006303172386, 0.000000000002666983691600854, 0.0000016330902276362, 0.000001664876689494737)

# Capture the current (pre-edit) B:F values for rows 2-10 before we overwrite anything,
# since rows 3-11 will receive the values currently in rows 2-10.
$oldValues = @{}
for ($r = 2; $r -le 10; $r++) {
    $rowVals = @()
    for ($c = 2; $c -le 6; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $oldValues[$r] = $rowVals
}

# Shift rows 10 -> 11 down to 2 -> 3 (process from bottom to top to avoid overwriting source data).
for ($r = 10; $r -ge 2; $r--) {
    $vals = $oldValues[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item($r + 1, $i + 2).Value2 = $vals[$i]
    }
}

# Insert the new row 2 values.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item(2, $i + 2).Value2 = $newRow2[$i]
}

# Column G (index 7) increments by 1 for every row 2..11.
for ($r = 2; $r -le 11; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 7).Value2 = $g + 1
}

$wb.Save()
